$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 44691
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 3500
$ws.Range("M2").Value = 3250
$ws.Range("N2").Value = "$/docena de matas"
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 542
$ws.Range("Q2").Value = 6

# Row 3 updates
$ws.Range("D3").Value = 44687
$ws.Range("J3").Value = 160

# Row 4 updates
$ws.Range("D4").Value = 44221
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 1300
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = 1420
$ws.Range("N4").Value = "$/atado"
$ws.Range("O4").Value = "Provincia de Diguillín"
$ws.Range("P4").Value = 1420
$ws.Range("Q4").Value = 1
